$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 364
$ws.Range("I9").Value = 72
$ws.Range("J9").Value = 480.8
$ws.Range("K9").Value = 72
$ws.Range("L9").Value = 480.8
$ws.Range("M9").Value = 97
$ws.Range("N9").Value = -818.8

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H132").Value = 822.6
$ws.Range("I132").Value = 822.6
$ws.Range("K132").Value = 2467.8
$ws.Range("M132").Value = 62.19999999999982

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 68999.5
$ws.Range("J92").Value = 68999.5
$ws.Range("L92").Value = 68999.5
$ws.Range("N92").Value = -73991.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3545.25
$ws.Range("I86").Value = 3454.3
$ws.Range("K86").Value = 3454.3
$ws.Range("M86").Value = -2331.3

$ws.Range("H89").Value = 3545.25
$ws.Range("I89").Value = 3454.3
$ws.Range("K89").Value = 17271.5
$ws.Range("M89").Value = -11655.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 350
$ws.Range("J16").Value = 400
$ws.Range("L16").Value = 400
$ws.Range("N16").Value = -974

$ws.Range("H22").Value = 4090.3157
$ws.Range("J22").Value = 2087.5
$ws.Range("L22").Value = 2087.5
$ws.Range("N22").Value = -2787.5

$ws.Range("H45").Value = 5450
$ws.Range("I45").Value = 4500
$ws.Range("J45").Value = 6400
$ws.Range("K45").Value = 4500
$ws.Range("L45").Value = 6400
$ws.Range("M45").Value = -3907
$ws.Range("N45").Value = -7586

$ws.Range("H58").Value = 1669
$ws.Range("I58").Value = 1503.5
$ws.Range("K58").Value = 1503.5
$ws.Range("M58").Value = -1300.5

$ws.Range("H88").Value = 18375
$ws.Range("J88").Value = 18500
$ws.Range("L88").Value = 18500
$ws.Range("N88").Value = -19312

$ws.Range("H91").Value = 18375
$ws.Range("J91").Value = 18500
$ws.Range("L91").Value = 18500
$ws.Range("N91").Value = -21308

$ws.Range("H92").Value = 49329.668
$ws.Range("J92").Value = 49329.668
$ws.Range("L92").Value = 49329.668
$ws.Range("N92").Value = -54321.668

$ws.Range("H113").Value = 350
$ws.Range("J113").Value = 400
$ws.Range("L113").Value = 400
$ws.Range("N113").Value = -4740

$ws.Range("H134").Value = 1719.3334
$ws.Range("I134").Value = 1615.5
$ws.Range("J134").Value = 1927
$ws.Range("K134").Value = 4846.5
$ws.Range("L134").Value = 5781
$ws.Range("M134").Value = -2311.5
$ws.Range("N134").Value = -10851

$ws.Range("H136").Value = 1669
$ws.Range("I136").Value = 1503.5
$ws.Range("K136").Value = 4510.5
$ws.Range("M136").Value = -1960.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 919.8
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H34").Value = 1204.7778
$ws.Range("J34").Value = 1568.1538
$ws.Range("L34").Value = 4704.4614
$ws.Range("N34").Value = -4872.4614

$ws.Range("H60").Value = 935.7143
$ws.Range("I60").Value = 935.7143
$ws.Range("K60").Value = 2807.1429
$ws.Range("M60").Value = -2556.1429

$ws.Range("H61").Value = 91.666664
$ws.Range("I61").Value = 91.666664
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 274.999992
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -59.99999200000002
$ws.Range("N61").ClearContents()

$ws.Range("H69").Value = 700
$ws.Range("J69").Value = 200
$ws.Range("L69").Value = 600
$ws.Range("N69").Value = -2222

$ws.Range("H72").Value = 700
$ws.Range("J72").Value = 200
$ws.Range("L72").Value = 1800
$ws.Range("N72").Value = -9912

$ws.Range("H109").Value = 2204.8572
$ws.Range("I109").Value = 72.333336
$ws.Range("K109").Value = 217.000008
$ws.Range("M109").Value = 822.999992

$ws.Range("H132").Value = 869.75
$ws.Range("I132").Value = 869.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7827.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5297.75
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 186.9
$ws.Range("I2").Value = 48.333332
$ws.Range("J2").Value = 300.27274
$ws.Range("K2").Value = 48.333332
$ws.Range("L2").Value = 300.27274
$ws.Range("M2").Value = 64.666668
$ws.Range("N2").Value = -526.27274

$ws.Range("H80").Value = 10249.25
$ws.Range("I80").Value = 5332.3335
$ws.Range("J80").Value = 25000
$ws.Range("K80").Value = 5332.3335
$ws.Range("L80").Value = 25000
$ws.Range("M80").Value = -4334.3335
$ws.Range("N80").Value = -26996

$ws.Range("H83").Value = 10249.25
$ws.Range("I83").Value = 5332.3335
$ws.Range("J83").Value = 25000
$ws.Range("K83").Value = 26661.6675
$ws.Range("L83").Value = 125000
$ws.Range("M83").Value = -21669.6675
$ws.Range("N83").Value = -134984

$ws.Range("H102").Value = 1487
$ws.Range("I102").Value = 1273.5
$ws.Range("J102").Value = 1914
$ws.Range("K102").Value = 1273.5
$ws.Range("L102").Value = 1914
$ws.Range("M102").Value = 348.5
$ws.Range("N102").Value = -5158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 540.6316
$ws.Range("I55").Value = 391.8889
$ws.Range("K55").Value = 391.8889
$ws.Range("M55").Value = -218.8889

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1798

$ws.Range("H68").Value = 5503
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5503
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5503
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7001

$ws.Range("H71").Value = 5503
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5503
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 27515
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -35003

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws.Range("H132").Value = 1499.25
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 1666
$ws.Range("K132").Value = 2997
$ws.Range("L132").Value = 4998
$ws.Range("M132").Value = -467
$ws.Range("N132").Value = -10058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29950
$ws.Range("J92").Value = 29950
$ws.Range("L92").Value = 29950
$ws.Range("N92").Value = -34942

$ws.Range("H107").Value = 1152.3636
$ws.Range("I107").Value = 1214
$ws.Range("J107").Value = 1078.4
$ws.Range("K107").Value = 3642
$ws.Range("L107").Value = 3235.2
$ws.Range("M107").Value = -1722
$ws.Range("N107").Value = -7075.200000000001

$ws.Range("H122").Value = 4996.5
$ws.Range("I122").Value = 4996.5
$ws.Range("K122").Value = 14989.5
$ws.Range("M122").Value = -12539.5

$ws.Range("H132").Value = 1730.909
$ws.Range("I132").Value = 1542.1428
$ws.Range("K132").Value = 4626.428400000001
$ws.Range("M132").Value = -2096.428400000001
